$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns whose contents move between rows (per the row-mapping below)
$cols = @("A", "B", "E", "F", "G", "H", "Q", "R", "AW", "AX")

# mapping: destination row -> source row (data from source row moves into destination row)
$mapping = @{
    3  = 4
    4  = 5
    5  = 6
    6  = 7
    7  = 9
    8  = 12
    9  = 14
    10 = 16
    11 = 17
    12 = 3
    13 = 8
    14 = 10
    15 = 11
    16 = 13
    17 = 15
}

# Snapshot all the "before" values for rows 3..17 for the tracked columns
$snapshot = @{}
foreach ($row in 3..17) {
    $rowVals = @{}
    foreach ($col in $cols) {
        $rowVals[$col] = $ws.Range("$col$row").Value2
    }
    $snapshot[$row] = $rowVals
}

# Apply the new values according to the mapping
foreach ($row in 3..17) {
    $srcRow = $mapping[$row]
    $srcVals = $snapshot[$srcRow]
    foreach ($col in $cols) {
        $ws.Range("$col$row").Value = $srcVals[$col]
    }
}

Write-Host "done"
